$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEAVE CREDITS")
$ws.Range("A1").Value = "test"
